$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 6, shifting rows 6-9 down to 7-10
$ws.Rows.Item(6).Insert()

# Copy the date cell style (s="2") from the row below (now row 7, formerly row 6)
$ws.Range("D7").Copy()
$ws.Range("D6").PasteSpecial(-4122) # xlPasteFormats

# Fill in the new row 6 values
$ws.Cells.Item(6, 1).Value = 1
$ws.Cells.Item(6, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(6, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(6, 4).Value = 44893
$ws.Cells.Item(6, 5).Value = 15
$ws.Cells.Item(6, 6).Value = "Fruta"
$ws.Cells.Item(6, 7).Value = 100108
$ws.Cells.Item(6, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(6, 9).Value = 100108007
$ws.Cells.Item(6, 10).Value = "Coco"
$ws.Cells.Item(6, 11).Value = "Sin especificar"
$ws.Cells.Item(6, 12).Value = "Primera"
$ws.Cells.Item(6, 13).Value = 80
$ws.Cells.Item(6, 14).Value = 21000
$ws.Cells.Item(6, 15).Value = 22000
$ws.Cells.Item(6, 16).Value = 21625
$ws.Cells.Item(6, 17).Value = "`$/malla 20 unidades"
$ws.Cells.Item(6, 18).Value = "Perú"
$ws.Cells.Item(6, 19).Value = 1081
$ws.Cells.Item(6, 20).Value = 20
